$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts name/type/major (and data) one column right
$ws.Range("A1").EntireColumn.Insert()

# Copy the header formatting (bold/border/centered) from the former first column header
# (now shifted to B1) onto the new header cell A1
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Populate the new "audit_id" column
$ws.Range("A1").Value = "audit_id"
$ws.Range("A2").Value = "is_0"
$ws.Range("A3").Value = "is_1"
$ws.Range("A4").Value = "cs_0"
$ws.Range("A5").Value = "cs_1"
$ws.Range("A6").Value = "ba_0"
$ws.Range("A7").Value = "ba_1"
$ws.Range("A8").Value = "bs_0"
$ws.Range("A9").Value = "bs_1"
